$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect numeric-looking "Price" text values (column D) from Excel
# auto-converting them to floating point numbers by switching the
# range to Text format before assigning, then restoring the default
# "Normal" style afterwards so the saved cells keep no style override,
# matching the original file.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# --- Price (column D) updates ---
$ws.Range("D2").Value = '26.884.09'
$ws.Range("D3").Value = '1.873.27'
$ws.Range("D5").Value = '301.95'
$ws.Range("D6").Value = '1.001'
$ws.Range("D7").Value = '0.5298'
$ws.Range("D8").Value = '0.3754'
$ws.Range("D10").Value = '21.59'
$ws.Range("D11").Value = '0.8848'
$ws.Range("D13").Value = '1.870.38'
$ws.Range("D14").Value = '92.98'
$ws.Range("D15").Value = '5.271'
$ws.Range("D16").Value = '1.002'
$ws.Range("D17").Value = '14.70'
$ws.Range("D18").Value = '0.000008537'
$ws.Range("D20").Value = '26.983.11'
$ws.Range("D21").Value = '4.971'
$ws.Range("D22").Value = '10.67'
$ws.Range("D23").Value = '6.376'
$ws.Range("D24").Value = '147.31'
$ws.Range("D25").Value = '2.261'
$ws.Range("D28").Value = '114.50'
$ws.Range("D29").Value = '4.737'
$ws.Range("D30").Value = '4.565'
$ws.Range("D31").Value = '0.09098'
$ws.Range("D32").Value = '0.7976'
$ws.Range("D33").Value = '0.04978'
$ws.Range("D34").Value = '1.172'
$ws.Range("D35").Value = '2.977'
$ws.Range("D36").Value = '3.200'
$ws.Range("D37").Value = '0.5844'
$ws.Range("D38").Value = '2.600'
$ws.Range("D39").Value = '1.072'
$ws.Range("D40").Value = '0.01948'
$ws.Range("D41").Value = '6.595'
$ws.Range("D42").Value = '8.878'
$ws.Range("D43").Value = '116.13'
$ws.Range("D44").Value = '0.5044'
$ws.Range("D45").Value = '0.1492'
$ws.Range("D47").Value = '9.956'
$ws.Range("D48").Value = '1.609'
$ws.Range("D49").Value = '37.97'
$ws.Range("D50").Value = '0.06028'
$ws.Range("D51").Value = '62.48'

$priceRange.Style = "Normal"

# --- Coin / Link / Volume(1h) updates ---
$ws.Range("E2").Value = '  -0.92%  '
$ws.Range("E3").Value = '  -1.39%  '
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("E7").Value = '  +1.30%  '
$ws.Range("E8").Value = '  -1.16%  '
$ws.Range("E9").Value = '  -1.71%  '
$ws.Range("E10").Value = '  +1.35%  '
$ws.Range("E11").Value = '  -2.27%  '
$ws.Range("E12").Value = '  -0.84%  '
$ws.Range("E13").Value = '  -1.35%  '
$ws.Range("E14").Value = '  -2.63%  '
$ws.Range("E15").Value = '  -1.46%  '
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("E17").Value = '  +0.20%  '
$ws.Range("E18").Value = '  -1.36%  '
$ws.Range("E19").Value = '  -0.23%  '
$ws.Range("E20").Value = '  -0.72%  '
$ws.Range("E21").Value = '  -2.91%  '
$ws.Range("E22").Value = '  -1.06%  '
$ws.Range("E23").Value = '  -1.41%  '
$ws.Range("E24").Value = '  -1.47%  '
$ws.Range("E25").Value = '  -2.90%  '
$ws.Range("E26").Value = '  -0.05%  '
$ws.Range("E27").Value = '  -1.55%  '
$ws.Range("E28").Value = '  -0.75%  '
$ws.Range("E29").Value = '  -1.79%  '
$ws.Range("E30").Value = '  -6.20%  '
$ws.Range("E31").Value = '  -1.47%  '
$ws.Range("E32").Value = '  +0.60%  '
$ws.Range("E33").Value = '  -1.36%  '
$ws.Range("E34").Value = '  -4.19%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("B36").Value = 'MXToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("E36").Value = '  -5.35%  '
$ws.Range("B37").Value = 'TheSandbox'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("E37").Value = '  +1.92%  '
$ws.Range("E38").Value = '  -1.90%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("E39").Value = '  -0.75%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("E40").Value = '  -2.28%  '
$ws.Range("E41").Value = '  -0.34%  '
$ws.Range("E42").Value = '  -1.58%  '
$ws.Range("E43").Value = '  -0.19%  '
$ws.Range("E44").Value = '  +2.98%  '
$ws.Range("E45").Value = '  -1.61%  '
$ws.Range("E46").Value = '  -0.33%  '
$ws.Range("E47").Value = '  -2.23%  '
$ws.Range("E48").Value = '  -1.89%  '
$ws.Range("E49").Value = '  -1.59%  '
$ws.Range("E50").Value = '  +1.22%  '
$ws.Range("E51").Value = '  -2.57%  '
